$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 226
$ws.Range("F3").Value = 54976
$ws.Range("F4").Value = 3118
$ws.Range("F5").Value = 5270
$ws.Range("F6").Value = 1190
$ws.Range("F8").Value = 861
$ws.Range("F9").Value = 27
$ws.Range("F10").Value = 1115
$ws.Range("F11").Value = 1376
$ws.Range("F12").Value = 125
$ws.Range("F14").Value = 224
$ws.Range("F15").Value = 397
$ws.Range("F16").Value = 62
$ws.Range("F19").Value = 73
$ws.Range("F20").Value = 67
$ws.Range("F21").Value = 5441
$ws.Range("F22").Value = 40
$ws.Range("F23").Value = 5329
$ws.Range("F24").Value = 9350
$ws.Range("F27").Value = 152
$ws.Range("F28").Value = 242
$ws.Range("F29").Value = 455
$ws.Range("F30").Value = 140
$ws.Range("F31").Value = 109
$ws.Range("F32").Value = 4278
$ws.Range("F33").Value = 294

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 141
$ws.Range("F6").Value = 10
$ws.Range("F10").Value = 1150
$ws.Range("F18").Value = 91

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 597

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 597
$ws.Range("F4").Value = 226
$ws.Range("F5").Value = 3118
$ws.Range("F7").Value = 141
$ws.Range("F8").Value = 1190
$ws.Range("F11").Value = 861
$ws.Range("F12").Value = 27
$ws.Range("F13").Value = 1115
$ws.Range("F14").Value = 10
$ws.Range("F15").Value = 1376
$ws.Range("F17").Value = 125
$ws.Range("F18").Value = 224
$ws.Range("F20").Value = 397
$ws.Range("F21").Value = 62
$ws.Range("F24").Value = 73
$ws.Range("F25").Value = 67
$ws.Range("F26").Value = 5441
$ws.Range("F27").Value = 40
$ws.Range("F28").Value = 5329
$ws.Range("F29").Value = 9350
$ws.Range("F33").Value = 152
$ws.Range("F34").Value = 242
$ws.Range("F35").Value = 455
$ws.Range("F39").Value = 140
$ws.Range("F40").Value = 109
$ws.Range("F41").Value = 4278
$ws.Range("F43").Value = 91
$ws.Range("F47").Value = 294
